$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = $null
$ws.Range("C2").Value = 5.3319794989134781
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = $null

$ws.Range("B3").Value = 5.6375100864256718
$ws.Range("C3").Value = 7.0164431192135419
$ws.Range("D3").Value = 8.7406576949142938
$ws.Range("E3").Value = 4.0263754227963036

$ws.Range("B1:E3").Select()
